$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the current row 26 (start of the "06/2025" block),
# shifting the existing data (June/May/April 2025 rows) down by 3 rows.
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(26).Insert()

# Fill in the new rows with the additional July 2025 (day 25-27) data.
$newData = @(
    @(25, 21050.1, 7, 2025, "07/2025"),
    @(26, 9960,    7, 2025, "07/2025"),
    @(27, 3677,    7, 2025, "07/2025")
)

$r = 26
foreach ($row in $newData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
Write-Host "done"
